# Commit: "update orientation for new amxEuler.m"
#
# 1) Bump the cached datetimeFigureOut field text (slide master + all
#    11 custom layouts) from 12/5/2019 -> 12/18/2019.
# 2) Re-orient the existing "Straight Arrow Connector 2" / "North" label
#    pair (used for the amxEuler.m heading diagram).
# 3) Add a second arrow + "East" label, duplicated from the existing
#    connector/label so they inherit the same line style/text formatting.

$p = $ppt.ActivePresentation

# --- 1. Date placeholder text (master + every custom layout) ---------------
$p.SlideMaster.Shapes.Item("Date Placeholder 3").TextFrame.TextRange.Text = "12/18/2019"

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $lyt = $p.SlideMaster.CustomLayouts.Item($li)
    for ($si = 1; $si -le $lyt.Shapes.Count; $si++) {
        $shp = $lyt.Shapes.Item($si)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "12/18/2019"
        }
    }
}

# --- 2. Work on slide 1 shapes ---------------------------------------------
$s = $p.Slides.Item(1)

$orientConnector = $s.Shapes.Item("Straight Arrow Connector 2")
$northLabel = $s.Shapes.Item("TextBox 24")

# Duplicate the connector + label *before* touching their geometry so the
# copies inherit the original "North" orientation/line style.
$newConnShapeRange = $orientConnector.Duplicate()
$newConnShape = $newConnShapeRange.Item(1)

$newLabelShapeRange = $northLabel.Duplicate()
$newLabelShape = $newLabelShapeRange.Item(1)

# --- Re-orient the original connector (North arrow) -------------------------
$orientConnector.Left = 807.2235107421875
$orientConnector.Top = 241.9252777099609375
$orientConnector.Width = 53.7174835205078125
$orientConnector.Height = 0.00007874015864217654
$orientConnector.HorizontalFlip = -1

# --- Move the "North" label to sit beside the re-oriented arrow -------------
$northLabel.Left = 749.63726806640625
$northLabel.Top = 229.808197021484375

# --- New "East" arrow (duplicate of the original connector) ----------------
$newConnShape.Name = "Straight Arrow Connector 27"
$newConnShape.Left = 861.88189697265625
$newConnShape.Top = 202.2111053466796875
$newConnShape.Width = 0
$newConnShape.Height = 40.834491729736328125

# --- New "East" label (duplicate of the "North" label) ----------------------
$newLabelShape.Name = "TextBox 29"
$newLabelShape.TextFrame.TextRange.Text = "East"
$newLabelShape.Left = 841.81097412109375
$newLabelShape.Top = 176.5753631591796875
$newLabelShape.Width = 38.26000213623046875
$newLabelShape.Height = 24.2344112396240234375
